$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.271.62'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '3.593.41'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.88'
$ws.Range('E5').Value = '  -2.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '190.97'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.634'
$ws.Range('E7').Value = '  -1.47%  '
$ws.Range('D8').Value = '3.586.45'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.179'
$ws.Range('E10').Value = '  -3.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.662'
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.65'
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000300'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.79'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = '4.167.27'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.19'
$ws.Range('E16').Value = '  +4.46%  '
$ws.Range('D17').Value = '3.582.43'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '70.171.45'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.55'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '19.64'
$ws.Range('E22').Value = '  +14.29%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '474.31'
$ws.Range('E23').Value = '  -5.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.13'
$ws.Range('E24').Value = '  -6.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.37'
$ws.Range('E25').Value = '  -1.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.70'
$ws.Range('E26').Value = '  -2.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.07'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.11'
$ws.Range('E28').Value = '  -0.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.29'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.80'
$ws.Range('E30').Value = '  +3.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.19'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.13'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '66.20'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '588.91'
$ws.Range('E35').Value = '  -3.83%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.56'
$ws.Range('E36').Value = '  +3.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '0.0₃0808'
$ws.Range('E38').Value = '  -3.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.400'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.54'
$ws.Range('E41').Value = '  -2.70%  '
$ws.Range('E42').Value = '  +7.15%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.17'
$ws.Range('E43').Value = '  +8.61%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '3.236.27'
$ws.Range('E44').Value = '  -2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.11'
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0447'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.63'
$ws.Range('E47').Value = '  +5.33%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.17'
$ws.Range('E51').Value = '  -2.30%  '
